$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 8082
$ws.Range("I40").Value = 3750
$ws.Range("J40").Value = 10248
$ws.Range("K40").Value = 3750
$ws.Range("L40").Value = 10248
$ws.Range("M40").Value = -3575
$ws.Range("N40").Value = -10598

$ws.Range("H103").Value = 472.4
$ws.Range("I103").Value = 210.875
$ws.Range("J103").Value = 771.2857
$ws.Range("K103").Value = 632.625
$ws.Range("L103").Value = 2313.8571
$ws.Range("M103").Value = -46.625
$ws.Range("N103").Value = -3485.8571

$ws.Range("H125").Value = 9228
$ws.Range("I125").Value = 9262.5
$ws.Range("K125").Value = 83362.5
$ws.Range("M125").Value = -80902.5

$ws.Range("H132").Value = 1923.129
$ws.Range("I132").Value = 1024.8846
$ws.Range("J132").Value = 6594
$ws.Range("K132").Value = 3074.6538
$ws.Range("L132").Value = 19782
$ws.Range("M132").Value = -544.6538
$ws.Range("N132").Value = -24842

$ws.Range("H137").Value = 4616.359
$ws.Range("I137").Value = 2804.577
$ws.Range("K137").Value = 8413.731
$ws.Range("M137").Value = -5863.731

$ws.Range("H138").Value = 5672.25
$ws.Range("I138").Value = 4617.7856
$ws.Range("J138").Value = 6106.4414
$ws.Range("K138").Value = 13853.3568
$ws.Range("L138").Value = 18319.3242
$ws.Range("M138").Value = -8713.356800000001
$ws.Range("N138").Value = -28599.3242

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 28000
$ws.Range("J24").Value = 28000
$ws.Range("L24").Value = 28000
$ws.Range("N24").Value = -28748

$ws.Range("H32").Value = 2482.3735
$ws.Range("I32").Value = 2158.3718
$ws.Range("K32").Value = 2158.3718
$ws.Range("M32").Value = -1871.3718

$ws.Range("H100").Value = 28000
$ws.Range("J100").Value = 28000
$ws.Range("L100").Value = 28000
$ws.Range("N100").Value = -30164

$ws.Range("H122").Value = 5158.2354
$ws.Range("I122").Value = 3355.375
$ws.Range("K122").Value = 10066.125
$ws.Range("M122").Value = -7616.125

$ws.Range("H132").Value = 8057.76
$ws.Range("I132").Value = 3034.5386
$ws.Range("J132").Value = 13499.583
$ws.Range("K132").Value = 9103.6158
$ws.Range("L132").Value = 40498.749
$ws.Range("M132").Value = -6573.6158
$ws.Range("N132").Value = -45558.749

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 89755.664
$ws.Range("J124").Value = 89755.664
$ws.Range("L124").Value = 89755.664
$ws.Range("N124").Value = -99575.664

$ws.Range("H134").Value = 2447.2307
$ws.Range("I134").Value = 1681.5758
$ws.Range("K134").Value = 5044.7274
$ws.Range("M134").Value = -2509.7274

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 315925.12
$ws.Range("I58").Value = 625992.4
$ws.Range("J58").Value = 5857.875
$ws.Range("K58").Value = 625992.4
$ws.Range("L58").Value = 5857.875
$ws.Range("M58").Value = -625789.4
$ws.Range("N58").Value = -6263.875

$ws.Range("H122").Value = 5152.3335
$ws.Range("I122").Value = 3900
$ws.Range("J122").Value = 5402.8
$ws.Range("K122").Value = 11700
$ws.Range("L122").Value = 16208.4
$ws.Range("M122").Value = -9250
$ws.Range("N122").Value = -21108.4

$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -54920

$ws.Range("H134").Value = 3252.3555
$ws.Range("I134").Value = 2362.0938
$ws.Range("K134").Value = 7086.2814
$ws.Range("M134").Value = -4551.2814

$ws.Range("H136").Value = 315925.12
$ws.Range("I136").Value = 625992.4
$ws.Range("J136").Value = 5857.875
$ws.Range("K136").Value = 1877977.2
$ws.Range("L136").Value = 17573.625
$ws.Range("M136").Value = -1875427.2
$ws.Range("N136").Value = -22673.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 48438.176
$ws.Range("I5").Value = 81029.10000000001
$ws.Range("J5").Value = 1879.7142
$ws.Range("K5").Value = 243087.3
$ws.Range("L5").Value = 5639.142599999999
$ws.Range("M5").Value = -242975.3
$ws.Range("N5").Value = -5863.142599999999

$ws.Range("H107").Value = 66724.31
$ws.Range("J107").Value = 75819.92999999999
$ws.Range("L107").Value = 227459.79
$ws.Range("N107").Value = -231299.79

$ws.Range("H132").Value = 4067.76
$ws.Range("I132").Value = 3542
$ws.Range("J132").Value = 4272.222
$ws.Range("K132").Value = 31878
$ws.Range("L132").Value = 38449.998
$ws.Range("M132").Value = -29348
$ws.Range("N132").Value = -43509.998

$ws.Range("H135").Value = 48438.176
$ws.Range("I135").Value = 81029.10000000001
$ws.Range("J135").Value = 1879.7142
$ws.Range("K135").Value = 729261.9
$ws.Range("L135").Value = 16917.4278
$ws.Range("M135").Value = -726726.9
$ws.Range("N135").Value = -21987.4278

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 239.25
$ws.Range("J2").Value = 66.666664
$ws.Range("L2").Value = 66.666664
$ws.Range("N2").Value = -292.666664

$ws.Range("H104").Value = 12671
$ws.Range("J104").Value = 12671
$ws.Range("L104").Value = 12671
$ws.Range("N104").Value = -19659

$ws.Range("H126").Value = 142860820
$ws.Range("J126").Value = 3937.3333
$ws.Range("L126").Value = 11811.9999
$ws.Range("N126").Value = -16751.9999

$ws.Range("H132").Value = 247153.58
$ws.Range("I132").Value = 273227
$ws.Range("K132").Value = 819681
$ws.Range("M132").Value = -817151

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 480240.1
$ws.Range("J40").Value = 6865.6665
$ws.Range("L40").Value = 6865.6665
$ws.Range("N40").Value = -7137.6665

$ws.Range("H122").Value = 2105191.8
$ws.Range("I122").Value = 2504977
$ws.Range("J122").Value = 1838668.4
$ws.Range("K122").Value = 7514931
$ws.Range("L122").Value = 5516005.199999999
$ws.Range("M122").Value = -7512481
$ws.Range("N122").Value = -5520905.199999999

$ws.Range("H135").Value = 67999.5
$ws.Range("J135").Value = 67999.5
$ws.Range("L135").Value = 67999.5
$ws.Range("N135").Value = -78139.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1611
$ws.Range("I100").Value = 222
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 444
$ws.Range("L100").Value = 6000
$ws.Range("M100").Value = 97
$ws.Range("N100").Value = -7082

$ws.Range("H132").Value = 3812.5957
$ws.Range("I132").Value = 3350.6
$ws.Range("J132").Value = 4627.8823
$ws.Range("K132").Value = 10051.8
$ws.Range("L132").Value = 13883.6469
$ws.Range("M132").Value = -7521.799999999999
$ws.Range("N132").Value = -18943.6469

$ws.Range("H136").Value = 2086.611
$ws.Range("I136").Value = 1977.4138
$ws.Range("K136").Value = 5932.2414
$ws.Range("M136").Value = -3382.2414
